$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 47706.715
$ws.Range("I33").Value = 55643.668
$ws.Range("J33").Value = 85
$ws.Range("K33").Value = 55643.668
$ws.Range("L33").Value = 85
$ws.Range("M33").Value = -55414.668
$ws.Range("N33").Value = -543

$ws.Range("H107").Value = 2500506.5
$ws.Range("I107").Value = 2500506.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2500506.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2498586.5
$ws.Range("N107").ClearContents()

$ws.Range("H111").Value = 713.9
$ws.Range("I111").Value = 682.1111
$ws.Range("J111").Value = 1000
$ws.Range("K111").Value = 2046.3333
$ws.Range("L111").Value = 3000
$ws.Range("M111").Value = 1020.6667
$ws.Range("N111").Value = -9134

$ws.Range("H125").Value = 2573.6365
$ws.Range("I125").Value = 860.7143
$ws.Range("J125").Value = 5571.25
$ws.Range("K125").Value = 7746.428699999999
$ws.Range("L125").Value = 50141.25
$ws.Range("M125").Value = -5286.428699999999
$ws.Range("N125").Value = -55061.25

$ws.Range("H132").Value = 1437.2439
$ws.Range("I132").Value = 1446.7949
$ws.Range("J132").Value = 1251
$ws.Range("K132").Value = 4340.384700000001
$ws.Range("L132").Value = 3753
$ws.Range("M132").Value = -1810.384700000001
$ws.Range("N132").Value = -8813

$ws.Range("H137").Value = 2410.182
$ws.Range("I137").Value = 2615.077
$ws.Range("J137").Value = 2277
$ws.Range("K137").Value = 7845.231000000001
$ws.Range("L137").Value = 6831
$ws.Range("M137").Value = -5295.231000000001
$ws.Range("N137").Value = -11931

$ws.Range("H138").Value = 4166.95
$ws.Range("I138").Value = 2227.1365
$ws.Range("J138").Value = 4714.077
$ws.Range("K138").Value = 6681.4095
$ws.Range("L138").Value = 14142.231
$ws.Range("M138").Value = -1541.4095
$ws.Range("N138").Value = -24422.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24101.508
$ws.Range("I32").Value = 21921.492
$ws.Range("J32").Value = 46991.668
$ws.Range("K32").Value = 21921.492
$ws.Range("L32").Value = 46991.668
$ws.Range("M32").Value = -21634.492
$ws.Range("N32").Value = -47565.668

$ws.Range("H61").Value = 1239.25
$ws.Range("I61").Value = 1105.2188
$ws.Range("J61").Value = 1596.6666
$ws.Range("K61").Value = 1105.2188
$ws.Range("L61").Value = 1596.6666
$ws.Range("M61").Value = -893.2188000000001
$ws.Range("N61").Value = -2020.6666

$ws.Range("H136").Value = 1239.25
$ws.Range("I136").Value = 1105.2188
$ws.Range("J136").Value = 1596.6666
$ws.Range("K136").Value = 3315.6564
$ws.Range("L136").Value = 4789.9998
$ws.Range("M136").Value = -765.6564000000003
$ws.Range("N136").Value = -9889.9998

$ws.Range("H138").Value = 42600
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 42600
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 42600
$ws.Range("N138").Value = -52880
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 45484950
$ws.Range("I20").Value = 65283.2
$ws.Range("J20").Value = 83334680
$ws.Range("K20").Value = 65283.2
$ws.Range("L20").Value = 83334680
$ws.Range("M20").Value = -65036.2
$ws.Range("N20").Value = -83335174

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1653.9385
$ws.Range("I31").Value = 707.5
$ws.Range("J31").Value = 2985.963
$ws.Range("K31").Value = 707.5
$ws.Range("L31").Value = 2985.963
$ws.Range("M31").Value = -412.5
$ws.Range("N31").Value = -3575.963

$ws.Range("H34").Value = 1653.9385
$ws.Range("I34").Value = 707.5
$ws.Range("J34").Value = 2985.963
$ws.Range("K34").Value = 707.5
$ws.Range("L34").Value = 2985.963
$ws.Range("M34").Value = -505.5
$ws.Range("N34").Value = -3389.963

$ws.Range("H134").Value = 2237.6
$ws.Range("I134").Value = 1618.875
$ws.Range("J134").Value = 4712.5
$ws.Range("K134").Value = 4856.625
$ws.Range("L134").Value = 14137.5
$ws.Range("M134").Value = -2321.625
$ws.Range("N134").Value = -19207.5

$ws.Range("H138").Value = 38425.25
$ws.Range("J138").Value = 38425.25
$ws.Range("L138").Value = 38425.25
$ws.Range("N138").Value = -48705.25

$ws.Range("H140").Value = 50523.727
$ws.Range("J140").Value = 50523.727
$ws.Range("L140").Value = 50523.727
$ws.Range("N140").Value = -60883.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1416.0212
$ws.Range("I68").Value = 734.6905
$ws.Range("J68").Value = 1966.3269
$ws.Range("K68").Value = 2204.0715
$ws.Range("L68").Value = 5898.9807
$ws.Range("M68").Value = -1393.0715
$ws.Range("N68").Value = -7520.9807

$ws.Range("H71").Value = 1416.0212
$ws.Range("I71").Value = 734.6905
$ws.Range("J71").Value = 1966.3269
$ws.Range("K71").Value = 6612.2145
$ws.Range("L71").Value = 17696.9421
$ws.Range("M71").Value = -2556.2145
$ws.Range("N71").Value = -25808.9421

$ws.Range("H107").Value = 371416.84
$ws.Range("I107").Value = 569.9286
$ws.Range("J107").Value = 1113110.8
$ws.Range("K107").Value = 1709.7858
$ws.Range("L107").Value = 3339332.4
$ws.Range("M107").Value = 210.2142000000001
$ws.Range("N107").Value = -3343172.4

$ws.Range("H131").Value = 25052960
$ws.Range("J131").Value = 3424.1765
$ws.Range("L131").Value = 10272.5295
$ws.Range("N131").Value = -20352.5295

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 488.5
$ws.Range("I107").Value = 470.42856
$ws.Range("J107").Value = 615
$ws.Range("K107").Value = 470.42856
$ws.Range("L107").Value = 615
$ws.Range("M107").Value = 1449.57144
$ws.Range("N107").Value = -4455

$ws.Range("H113").Value = 1444.5555
$ws.Range("I113").Value = 1005.8571
$ws.Range("K113").Value = 1005.8571
$ws.Range("M113").Value = 1164.1429

$ws.Range("H126").Value = 2105.875
$ws.Range("I126").Value = 1898.3334
$ws.Range("J126").Value = 2728.5
$ws.Range("K126").Value = 5695.0002
$ws.Range("L126").Value = 8185.5
$ws.Range("M126").Value = -3225.0002
$ws.Range("N126").Value = -13125.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2215.125
$ws.Range("I16").Value = 2215.125
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2215.125
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2045.125
$ws.Range("N16").ClearContents()

$ws.Range("H40").Value = 3667.88
$ws.Range("I40").Value = 3826.9092
$ws.Range("J40").Value = 2501.6667
$ws.Range("K40").Value = 3826.9092
$ws.Range("L40").Value = 2501.6667
$ws.Range("M40").Value = -3690.9092
$ws.Range("N40").Value = -2773.6667

$ws.Range("H139").Value = 57450.25
$ws.Range("J139").Value = 57450.25
$ws.Range("L139").Value = 57450.25
$ws.Range("N139").Value = -67730.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 321.1
$ws.Range("I113").Value = 348.45456
$ws.Range("K113").Value = 1045.36368
$ws.Range("M113").Value = 1124.63632
